$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 2, column B: rename RLIItemRateCountList -> RLIdValueList
$ws.Range("B2").Value = "RLIdValueList"

# Delete rows 8 through 14 (old data no longer needed), keep rows 1-7 with new content
$ws.Range("A8:B14").EntireRow.Delete()

# Rewrite rows 4-7 with the rearranged data (write the short strings first so the
# shared-strings table ends up ordered the same way as the target workbook)
$ws.Range("A5").Value = 22031201
$ws.Range("B5").Value = "22035001;5"

$ws.Range("A6").Value = 22031202
$ws.Range("B6").Value = "22035002;5"

$ws.Range("A7").Value = 22031203
$ws.Range("B7").Value = "22035003;5"

$ws.Range("A4").Value = 22031001
$ws.Range("B4").Value = "22033001;5|22033002;5|22032007;1|22033013;5|22033014;3|22033015;3"

# Resize the table to the new data extent
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B7"))

# Update the workbook theme's Background 1 (lt1) color, as was done via
# Page Layout > Colors > Customize Colors in the original edit
$tcs = $wb.Theme.ThemeColorScheme
$tcs.Colors(2).RGB = 13494986

# Move selection to match target view state
$ws.Range("B4").Select()
